# Corrected Calibration and Ingest Sheets for Coastal Gliders
# Changed FLORT cal values for angular resolution to 1.076, changed
# Scattering Angle to 124.

$wb = $excel.ActiveWorkbook
$wsCal = $wb.Worksheets.Item("Asset_Cal_Info")
$wsMoorings = $wb.Worksheets.Item("Moorings")

# CC_scattering_angle (row 4): 117 -> 124
$wsCal.Range("F4").Value = 124

# CC_angular_resolution (row 6): 1.08 -> 1.076
$wsCal.Range("F6").Value = 1.076

# Reflect the author's last selected cell on Asset_Cal_Info (D22 -> E22)
# without changing which sheet tab is active (Moorings stays active).
$wsCal.Range("E22").Select()
$wsMoorings.Activate()
